$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 416838.6
$ws.Range("I11").Value = 416838.6
$ws.Range("K11").Value = 416838.6
$ws.Range("M11").Value = -416698.6

$ws.Range("H62").Value = 7041.933
$ws.Range("I62").Value = 5785
$ws.Range("K62").Value = 5785
$ws.Range("M62").Value = -5161

$ws.Range("H65").Value = 7041.933
$ws.Range("I65").Value = 5785
$ws.Range("K65").Value = 28925
$ws.Range("M65").Value = -25805

$ws.Range("H74").Value = 13830.526
$ws.Range("I74").Value = 14518.8
$ws.Range("K74").Value = 14518.8
$ws.Range("M74").Value = -13582.8

$ws.Range("H77").Value = 13830.526
$ws.Range("I77").Value = 14518.8
$ws.Range("K77").Value = 72594
$ws.Range("M77").Value = -67914

$ws.Range("H138").Value = 2947.4912
$ws.Range("J138").Value = 3178.6667
$ws.Range("L138").Value = 9536.000100000001
$ws.Range("N138").Value = -19816.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 701.5333000000001
$ws.Range("I102").Value = 785.3333
$ws.Range("J102").Value = 366.33334
$ws.Range("K102").Value = 785.3333
$ws.Range("L102").Value = 366.33334
$ws.Range("M102").Value = 836.6667
$ws.Range("N102").Value = -3610.33334

$ws.Range("H110").Value = 1987.5834
$ws.Range("I110").Value = 1987.5834
$ws.Range("K110").Value = 1987.5834
$ws.Range("M110").Value = 57.41660000000002

$ws.Range("H111").Value = 12500
$ws.Range("J111").Value = 12500
$ws.Range("L111").Value = 12500
$ws.Range("N111").Value = -20680

$ws.Range("H122").Value = 2501.5356
$ws.Range("I122").Value = 2116.4583
$ws.Range("K122").Value = 6349.374899999999
$ws.Range("M122").Value = -3899.374899999999

$ws.Range("H132").Value = 2775.84
$ws.Range("I132").Value = 2156.7646
$ws.Range("K132").Value = 6470.293799999999
$ws.Range("M132").Value = -3940.293799999999

$ws.Range("H140").Value = 15300
$ws.Range("J140").Value = 15300
$ws.Range("L140").Value = 15300
$ws.Range("N140").Value = -25660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 17499.5
$ws.Range("I26").Value = 17499.5
$ws.Range("K26").Value = 17499.5
$ws.Range("M26").Value = -17207.5

$ws.Range("H64").Value = 460.25
$ws.Range("I64").Value = 586.1429000000001
$ws.Range("K64").Value = 586.1429000000001
$ws.Range("M64").Value = -361.1429000000001

$ws.Range("H67").Value = 460.25
$ws.Range("I67").Value = 586.1429000000001
$ws.Range("K67").Value = 586.1429000000001
$ws.Range("M67").Value = 193.8570999999999

$ws.Range("H86").Value = 2429.1482
$ws.Range("I86").Value = 1672.9474
$ws.Range("J86").Value = 4225.125
$ws.Range("K86").Value = 1672.9474
$ws.Range("L86").Value = 4225.125
$ws.Range("M86").Value = -549.9474
$ws.Range("N86").Value = -6471.125

$ws.Range("H89").Value = 2429.1482
$ws.Range("I89").Value = 1672.9474
$ws.Range("J89").Value = 4225.125
$ws.Range("K89").Value = 8364.737000000001
$ws.Range("L89").Value = 21125.625
$ws.Range("M89").Value = -2748.737000000001
$ws.Range("N89").Value = -32357.625

$ws.Range("H99").Value = 1910.1
$ws.Range("J99").Value = 1013.5
$ws.Range("L99").Value = 1013.5
$ws.Range("N99").Value = -4009.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1199
$ws.Range("I16").Value = 1199
$ws.Range("K16").Value = 1199
$ws.Range("M16").Value = -912

$ws.Range("H31").Value = 13410.611
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 13410.611
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H99").Value = 2471.8823
$ws.Range("I99").Value = 1927.4706
$ws.Range("J99").Value = 3016.2942
$ws.Range("K99").Value = 1927.4706
$ws.Range("L99").Value = 3016.2942
$ws.Range("M99").Value = -429.4706000000001
$ws.Range("N99").Value = -6012.2942

$ws.Range("H113").Value = 1199
$ws.Range("I113").Value = 1199
$ws.Range("K113").Value = 1199
$ws.Range("M113").Value = 971

$ws.Range("H126").Value = 2471.8823
$ws.Range("I126").Value = 1927.4706
$ws.Range("J126").Value = 3016.2942
$ws.Range("K126").Value = 5782.4118
$ws.Range("L126").Value = 9048.882599999999
$ws.Range("M126").Value = -3312.4118
$ws.Range("N126").Value = -13988.8826

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 74.5
$ws.Range("J12").Value = 88.125
$ws.Range("L12").Value = 264.375
$ws.Range("N12").Value = -610.375

$ws.Range("H61").Value = 158.77777
$ws.Range("I61").Value = 173.625
$ws.Range("K61").Value = 520.875
$ws.Range("M61").Value = -305.875

$ws.Range("H121").Value = 1789.25
$ws.Range("I121").Value = 1859.6364
$ws.Range("J121").Value = 1703.2222
$ws.Range("K121").Value = 5578.9092
$ws.Range("L121").Value = 5109.6666
$ws.Range("M121").Value = -4268.9092
$ws.Range("N121").Value = -7729.6666

$ws.Range("H129").Value = 3908.6667
$ws.Range("I129").Value = 1510.8889
$ws.Range("J129").Value = 7505.3335
$ws.Range("K129").Value = 4532.6667
$ws.Range("L129").Value = 22516.0005
$ws.Range("M129").Value = 467.3333000000002
$ws.Range("N129").Value = -32516.0005

$ws.Range("H137").Value = 4750.7856
$ws.Range("I137").Value = 2566.1667
$ws.Range("J137").Value = 6389.25
$ws.Range("K137").Value = 7698.500100000001
$ws.Range("L137").Value = 19167.75
$ws.Range("M137").Value = -2598.500100000001
$ws.Range("N137").Value = -29367.75

$ws.Range("H139").Value = 2290.7058
$ws.Range("I139").Value = 2277.625
$ws.Range("K139").Value = 6832.875
$ws.Range("M139").Value = -1692.875

$ws.Range("H140").Value = 1570.24
$ws.Range("I140").Value = 1055.6
$ws.Range("J140").Value = 1913.3334
$ws.Range("K140").Value = 3166.8
$ws.Range("L140").Value = 5740.0002
$ws.Range("M140").Value = 2013.2
$ws.Range("N140").Value = -16100.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1069.3334
$ws.Range("I97").Value = 983.4
$ws.Range("K97").Value = 983.4
$ws.Range("M97").Value = -487.4

$ws.Range("H113").Value = 6813
$ws.Range("I113").Value = 2713
$ws.Range("J113").Value = 51913
$ws.Range("K113").Value = 2713
$ws.Range("L113").Value = 51913
$ws.Range("M113").Value = -543
$ws.Range("N113").Value = -56253

$ws.Range("H132").Value = 4035.1082
$ws.Range("J132").Value = 11523.4
$ws.Range("L132").Value = 34570.2
$ws.Range("N132").Value = -39630.2

$ws.Range("H136").Value = 21310.814
$ws.Range("J136").Value = 21310.814
$ws.Range("L136").Value = 63932.442
$ws.Range("N136").Value = -69032.442

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 204158
$ws.Range("I7").Value = 204158
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 204158
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -204046
$ws.Range("N7").ClearContents()

$ws.Range("H40").Value = 7039.8335
$ws.Range("I40").Value = 3881
$ws.Range("J40").Value = 8254.77
$ws.Range("K40").Value = 3881
$ws.Range("L40").Value = 8254.77
$ws.Range("M40").Value = -3745
$ws.Range("N40").Value = -8526.77

$ws.Range("H93").Value = 1725.2273
$ws.Range("I93").Value = 1137.0667
$ws.Range("J93").Value = 2985.5715
$ws.Range("K93").Value = 1137.0667
$ws.Range("L93").Value = 2985.5715
$ws.Range("M93").Value = 110.9332999999999
$ws.Range("N93").Value = -5481.5715

$ws.Range("H126").Value = 204158
$ws.Range("I126").Value = 204158
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 612474
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -610004
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 29970
$ws.Range("J69").Value = 29970
$ws.Range("L69").Value = 29970
$ws.Range("N69").Value = -31468

$ws.Range("H72").Value = 29970
$ws.Range("J72").Value = 29970
$ws.Range("L72").Value = 89910
$ws.Range("N72").Value = -97398

$ws.Range("H126").Value = 1251.3334
$ws.Range("I126").Value = 1127
$ws.Range("K126").Value = 3381
$ws.Range("M126").Value = -911

$ws.Range("H132").Value = 3114.84
$ws.Range("I132").Value = 2505.0417
$ws.Range("K132").Value = 7515.125100000001
$ws.Range("M132").Value = -4985.125100000001
